$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documents")

$url = "https://dubaiholding-my.sharepoint.com/:b:/g/personal/arun_naidu_dhre_ae/EUysg7T3ajdIsj6J81h2b_8BqMv40tjbbnmHigxS5cyl8w?email=Arun.Naidu%40dhre.ae&e=fNZcxb"

# Add the new "PC#02" row of data (row 9), mirroring the existing row 8 (PC#01)
$ws.Range("A9").Value = "PC2.0"
$ws.Range("B9").Value = "DMS 149600"
$ws.Range("C9").Value = "PC#02 - Payment Certificate.pdf"
$ws.Range("D9").Value = "PAYMENT CERTIFICATE"
$ws.Range("E9").Value = 45828
$ws.Range("F9").Value = "PC # 02.0"
$ws.Range("G9").Value = "Payment Certificate # 02"
$ws.Range("H9").Value = "PC-02.PDF"
$ws.Range("I9").Value = $url

# Create the hyperlink for I9, then copy I8's formatting (Hyperlink cell style) onto it
$ws.Hyperlinks.Add($ws.Range("I9"), $url)
$ws.Range("I8").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = $url

$ws.Range("J9").Value = 45828
$ws.Range("K9").Value = 0

# Switch the active/selected sheet from Activities to Documents
$ws.Activate()
$ws.Range("K9").Select()
